$d = $word.ActiveDocument

$pairs = @(
    @("871÷3=", "396÷4="),
    @("773÷6=", "944÷3="),
    @("693÷3=", "930÷9="),
    @("213÷2=", "147÷3="),
    @("358÷8=", "962÷6="),
    @("973÷9=", "671÷8="),
    @("293÷8=", "906÷3="),
    @("828÷2=", "555÷5="),
    @("386÷5=", "316÷2="),
    @("403÷5=", "556÷6="),
    @("132÷6=", "185÷6="),
    @("203÷9=", "504÷8="),
    @("986÷7=", "375÷2="),
    @("445÷6=", "422÷5="),
    @("841÷3=", "526÷7="),
    @("935÷3=", "254÷2="),
    @("235÷2=", "490÷4="),
    @("208÷5=", "514÷2="),
    @("689÷6=", "199÷2="),
    @("598÷7=", "200÷2="),
    @("899÷2=", "531÷5="),
    @("604÷4=", "152÷7="),
    @("711÷8=", "670÷8="),
    @("154÷7=", "545÷4="),
    @("424÷4=", "723÷2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
